# Weekly refresh of the Fruta / hortaliza "Vega Monumental Concepcion - Tuna" subconjunto.
# The price series got reshuffled (dates/qualities/volumes/prices moved between rows) and
# one more historical record (the old row 12) was appended as a new row 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: append a new row 16 that duplicates the ORIGINAL row 12 (values + date style),
# before row 12 itself gets overwritten with its new values below.
$ws.Range("A12:T12").Copy() | Out-Null
$ws.Range("A16:T16").PasteSpecial() | Out-Null
$ws.Cells.Item(16, 4).NumberFormat = $ws.Cells.Item(12, 4).NumberFormat

# Step 2: columns -> A=1 B=2 C=3 D=4(Fecha) ... L=12(Calidad) M=13(Volumen)
#          N=14(Precio min) O=15(Precio max) P=16(Precio prom) Q=17(Unidad)
#          R=18(Origen) S=19(Precio $/Kg) T=20(Kg/unidad)
# Apply the per-cell value updates for rows 2-15.
$ws.Cells.Item(2, 4).Value = 45014
$ws.Cells.Item(2, 13).Value = 50
$ws.Cells.Item(2, 14).Value = 13000
$ws.Cells.Item(2, 15).Value = 14000
$ws.Cells.Item(2, 16).Value = 13600
$ws.Cells.Item(2, 19).Value = 756
$ws.Cells.Item(3, 4).Value = 45014
$ws.Cells.Item(3, 12).Value = "Segunda"
$ws.Cells.Item(3, 13).Value = 20
$ws.Cells.Item(3, 14).Value = 10000
$ws.Cells.Item(3, 15).Value = 10000
$ws.Cells.Item(3, 16).Value = 10000
$ws.Cells.Item(3, 19).Value = 556
$ws.Cells.Item(4, 4).Value = 45155
$ws.Cells.Item(4, 13).Value = 40
$ws.Cells.Item(4, 14).Value = 25000
$ws.Cells.Item(4, 15).Value = 26000
$ws.Cells.Item(4, 16).Value = 25500
$ws.Cells.Item(4, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(4, 19).Value = 1417
$ws.Cells.Item(5, 4).Value = 44819
$ws.Cells.Item(5, 13).Value = 100
$ws.Cells.Item(5, 14).Value = 25000
$ws.Cells.Item(5, 15).Value = 26000
$ws.Cells.Item(5, 16).Value = 25500
$ws.Cells.Item(5, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(5, 19).Value = 1417
$ws.Cells.Item(6, 4).Value = 44280
$ws.Cells.Item(6, 14).Value = 14000
$ws.Cells.Item(6, 15).Value = 15000
$ws.Cells.Item(6, 16).Value = 14500
$ws.Cells.Item(6, 19).Value = 806
$ws.Cells.Item(7, 12).Value = "Segunda"
$ws.Cells.Item(7, 13).Value = 50
$ws.Cells.Item(7, 14).Value = 12000
$ws.Cells.Item(7, 15).Value = 12000
$ws.Cells.Item(7, 16).Value = 12000
$ws.Cells.Item(7, 19).Value = 667
$ws.Cells.Item(8, 4).Value = 45002
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 100
$ws.Cells.Item(8, 15).Value = 13000
$ws.Cells.Item(8, 16).Value = 12500
$ws.Cells.Item(8, 19).Value = 694
$ws.Cells.Item(10, 4).Value = 44316
$ws.Cells.Item(10, 13).Value = 50
$ws.Cells.Item(10, 15).Value = 20000
$ws.Cells.Item(10, 16).Value = 20000
$ws.Cells.Item(10, 19).Value = 1111
$ws.Cells.Item(11, 4).Value = 45044
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 100
$ws.Cells.Item(11, 14).Value = 17000
$ws.Cells.Item(11, 16).Value = 17500
$ws.Cells.Item(11, 19).Value = 972
$ws.Cells.Item(12, 4).Value = 45030
$ws.Cells.Item(12, 14).Value = 15000
$ws.Cells.Item(12, 15).Value = 16000
$ws.Cells.Item(12, 16).Value = 15500
$ws.Cells.Item(12, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(12, 19).Value = 861
$ws.Cells.Item(13, 4).Value = 44687
$ws.Cells.Item(13, 13).Value = 100
$ws.Cells.Item(13, 14).Value = 18000
$ws.Cells.Item(13, 15).Value = 19000
$ws.Cells.Item(13, 16).Value = 18500
$ws.Cells.Item(13, 19).Value = 1028
$ws.Cells.Item(14, 4).Value = 44699
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 100
$ws.Cells.Item(14, 14).Value = 20000
$ws.Cells.Item(14, 15).Value = 22000
$ws.Cells.Item(14, 16).Value = 21000
$ws.Cells.Item(14, 19).Value = 1167
$ws.Cells.Item(15, 4).Value = 44699
$ws.Cells.Item(15, 12).Value = "Segunda"
$ws.Cells.Item(15, 13).Value = 50
$ws.Cells.Item(15, 14).Value = 18000
$ws.Cells.Item(15, 15).Value = 18000
$ws.Cells.Item(15, 16).Value = 18000
$ws.Cells.Item(15, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(15, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(15, 19).Value = 1000
$ws.Cells.Item(15, 20).Value = 18
